$wb = $excel.ActiveWorkbook

# --- Rename the original sheet ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "LAB 3 DigitalMarketingData"

# --- Add the new sheet right after it ---
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "LAB 3 Assignments "

# --- Update view/selection on the data sheet ---
$ws1.Activate()
$ws1.Range("A1:K11").Select()
$excel.ActiveWindow.Zoom = 85

# --- Populate the assignment sheet ---
$ws2.Columns.Item(2).ColumnWidth = 124
$ws2.Range("B3").Value = "📘 PART 1 — DATA ENTRY (Already Completed in Your File)`nYour dataset contains:`nFacebook Ad Spend (NPR)`nWebsite Visits`nLeads`nOrders`nConversion calculations`nCost calculations"
$ws2.Rows.Item([int]"3").RowHeight = 297
$ws2.Range("B3").WrapText = $true
$ws2.Range("B3").VerticalAlignment = -4160
$ws2.Range("B3").Copy()
$ws2.Range("B10").Value = "## 📌 PART 2 — REQUIRED EXCEL FORMULAS`nEnter these formulas EXACTLY:`n### 1. Visit → Lead Conversion`n= D2 / C2`n### 2. Lead → Order Conversion`n= E2 / D2`n### 3. Visit → Order Conversion`n= E2 / C2`n### 4. Cost per Visit`n= B2 / C2`n### 5. Cost per Lead`n= B2 / D2`n### 6. Cost per Order`n(Most important metric)`n= B2 / E2`n### 7. Average Orders`n=AVERAGE(E2:E11)`n### 8. Correlation (Relationship Strength)`nAd Spend vs Orders`n=CORREL(B2:B11, E2:E11)`nWebsite Visits vs Orders`n=CORREL(C2:C11, E2:E11)`nLeads vs Orders`n=CORREL(D2:D11, E2:E11)"
$ws2.Rows.Item([int]"10").RowHeight = 409.5
$ws2.Range("B10").PasteSpecial(-4122)
$ws2.Range("B20").Value = "## 📊 PART 3 — VISUALIZATION TASKS`n### 📌 Chart 1 — Line Chart (Trend Over Time)`nColumns to highlight:`nDay`nWebsite Visits`nLeads`nOrders`nInsert → Line Chart`nGive title: Daily Marketing Performance Trend`n### 📌 Chart 2 — Column Chart (Daily Orders)`nHighlight:`nDay`nOrders`nInsert → Column Chart`nTitle: New Customers Per Day`n### 📌 Chart 3 — Scatter Plot (Ad Spend → Orders)`nHighlight:`nAd Spend`nOrders`nInsert → Scatter Chart`nRight-click → Add Trendline`nCheck both:`n✔ Show Equation`n✔ Display R² Value`nInterpretation Expected:`nDoes increasing ad spend increase orders?`n### 📌 Chart 4 — Funnel Chart (Visitors → Leads → Orders)`nCreate this table:`nStage	Count`nWebsite Visits	=SUM(C2:C11)`nLeads	=SUM(D2:D11)`nOrders	=SUM(E2:E11)`nHighlight → Insert → Funnel`n## 🔍 PART 4 — ANALYSIS QUESTIONS (Students Must Answer)`nQ1. What is the average Visit → Order conversion rate?`n(Use formula results)`nQ2. Which metric has the strongest correlation with Orders?`n(Interpret CORREL results)`nQ3. Is Cost per Order increasing, decreasing, or stable?`n(Observe Column K)`nQ4. Which day had the best marketing performance? Why?"
$ws2.Rows.Item([int]"20").RowHeight = 409.5
$ws2.Range("B20").PasteSpecial(-4122)
$ws2.Range("B21").PasteSpecial(-4122)
$ws2.Range("B40").Value = "## 📘 PART 5 — SWOT ANALYSIS (Based on Your Excel Output)`nStudents must fill in:`n🟩 Strengths (Use correlation, low CPO, lead trends)`nExample: Strong positive relation between spend and orders.`n🟨 Weaknesses (Low conversion rate, rising CPO)`n🟦 Opportunities (Growing digital adoption, retargeting potential)`n🟥 Threats (Increasing CPC/competition)"
$ws2.Rows.Item([int]"40").RowHeight = 165
$ws2.Range("B40").PasteSpecial(-4122)
$ws2.Range("B63").Value = "## 🔮 PART 6 — FORECASTING IN EXCEL`n### 📌 Task 1 — Forecast Orders for Day 11`nUse:`n=FORECAST.LINEAR(11, E2:E11, A2:A11)`n### 📌 Task 2 — Forecast Orders for Ad Spend NPR 7000`n=FORECAST.LINEAR(7000, E2:E11, B2:B11)`n### 📌 Task 3 — Interpret Trendline Equation`nFrom Scatter → Trendline`nYou will see something like:`nOrders = 0.0048 * AdSpend + 1.5`nStudents must answer:`nWhat does the slope mean?`nIf SastoBazar spends 10,000 NPR, how many orders will they get?`nFormula: =0.0048 * 10000 + 1.5"
$ws2.Rows.Item([int]"63").RowHeight = 270
$ws2.Range("B63").PasteSpecial(-4122)
$ws2.Range("B89").Value = "s"
$ws2.Range("B89").PasteSpecial(-4122)

# --- Final view state: Assignments sheet active ---
$ws2.Range("A89").Select()
$ws2.Activate()
$excel.CutCopyMode = $false
